$d = $word.ActiveDocument

# The last paragraph currently holds the (hidden) _GoBack bookmark somewhere
# in the middle of its text. We need to end up with three brand-new plain
# paragraphs appended after it, and the _GoBack bookmark relocated to the
# very end of the text of the last of those new paragraphs.

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# Drop the existing bookmark - it will be re-created at the new location
# once all of the new content is in place.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Split off a new paragraph right after the last paragraph's text (before
# its paragraph mark).
$insertPos = $lastPara.Range.End - 1
$r = $d.Range($insertPos, $insertPos)
$r.InsertParagraphAfter()

# --- New paragraph 1: "UWAGA: ..." ---
$p1 = $d.Paragraphs.Item($lastParaIndex + 1)
$p1.Style = "Normal"
$p1.Range.Text = "UWAGA:  a może gdyby się udało zmniejszyć rozmiar m macierzy (liczbę wierszy, kosztem liczby kolumn oczywiście) to uzyskalibyśmy przyspieszenie? Wtedy jeden woksel nie mieściłby się w jednej kolumnie. Czy takie rozwiązanie jest możliwe? – zapytać się Gonzalo."
$p1.Range.InsertParagraphAfter()

# --- New paragraph 2: "- przed tym ..." ---
$p2 = $d.Paragraphs.Item($lastParaIndex + 2)
$p2.Style = "Normal"
$p2.Range.Text = "- przed tym jak spróbujesz ich zapytać, wykonaj w ogóle testy (nie zwracając większej uwagi na to czy taka operacja ma sens w świetle naszego algorytmu) czy to poprawia wydajność programu!"
$p2.Range.InsertParagraphAfter()

# --- New paragraph 3: "- zapytać ich ..." (with trailing placeholder char) ---
$p3 = $d.Paragraphs.Item($lastParaIndex + 3)
$p3.Style = "Normal"
$p3.Range.Text = "- zapytać ich jak im się współpracuje z lekarzami w hiszpani/ WalencjiX"

# Add the _GoBack bookmark right before the placeholder char 'X' (a safe,
# non-boundary position), then remove the placeholder so the bookmark ends
# up collapsed at the true end of the paragraph's text.
$bmPos = $p3.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($p3.Range.End - 2, $p3.Range.End - 1)
$placeholderRange.Delete()
